# "new wat script implementation"
# Adds a new test-case row (WAT29 / WAT-199) to the bottom of the
# "Test Cases" sheet, right after the existing last row (row 36),
# and moves the active selection to C14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Duplicate the formatting (borders/wrap/etc.) of the last existing row
# (row 36) into the new row 37 - this is how the row was authored
# originally (copy last row, then overwrite the text), and it keeps the
# new cells' styles in line with the rest of the sheet instead of
# minting brand-new style entries.
$ws.Range("A36:E36").Copy($ws.Range("A37:E37"))

# Overwrite with the new test case's data. Column order here controls
# the order new shared strings get appended in: TCID, Description, then
# JIRA ID - matching WAT29 / description / WAT-199.
$ws.Range("A37").Value() = "WAT29"
$ws.Range("C37").Value() = "Verify that System must display ""Top Journals"" section listing a maximum of three journal titles."
$ws.Range("B37").Value() = "WAT-199"
# D37 (Runmode) and E37 (Results) already carry over "Y" / blank from
# the copied row 36, matching the target row.

# Move the active selection, as recorded in the sheet view.
$ws.Range("C14").Select() | Out-Null
